# Update TPM-derived NATMI ligand-receptor metrics for Fgf1-Fgfr3 sheet
# per the new TPM re-run (commit: "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.729797666666667
$ws.Range("H2").Value = 5.189393000000001
$ws.Range("I2").Value = 0.06436583050179444
$ws.Range("J2").Value = 0.06436583050179444
$ws.Range("M2").Value = 4.959409333333333
$ws.Range("N2").Value = 14.878228
$ws.Range("O2").Value = 0.8271666313262851
$ws.Range("P2").Value = 0.8271666313262852
$ws.Range("Q2").Value = 8.57877469284489
$ws.Range("R2").Value = 77.20897223560401
$ws.Range("S2").Value = 0.05324126718868796
$ws.Range("T2").Value = 0.05324126718868796
$ws.Range("G3").Value = 1.729797666666667
$ws.Range("H3").Value = 5.189393000000001
$ws.Range("I3").Value = 0.06436583050179444
$ws.Range("J3").Value = 0.06436583050179444
$ws.Range("O3").Value = 0.09421438109281059
$ws.Range("P3").Value = 0.09421438109281059
$ws.Range("Q3").Value = 0.977123493152889
$ws.Range("R3").Value = 8.794111438376001
$ws.Range("S3").Value = 0.006064186884251313
$ws.Range("T3").Value = 0.006064186884251313
$ws.Range("G4").Value = 1.729797666666667
$ws.Range("H4").Value = 5.189393000000001
$ws.Range("I4").Value = 0.06436583050179444
$ws.Range("J4").Value = 0.06436583050179444
$ws.Range("O4").Value = 0.07861898758090437
$ws.Range("P4").Value = 0.07861898758090438
$ws.Range("Q4").Value = 0.8153793389304446
$ws.Range("R4").Value = 7.338414050374001
$ws.Range("S4").Value = 0.005060376428855173
$ws.Range("T4").Value = 0.005060376428855174
$ws.Range("I5").Value = 0.2200595722726403
$ws.Range("J5").Value = 0.2200595722726403
$ws.Range("M5").Value = 4.959409333333333
$ws.Range("N5").Value = 14.878228
$ws.Range("O5").Value = 0.8271666313262851
$ws.Range("P5").Value = 0.8271666313262852
$ws.Range("Q5").Value = 29.32987075305689
$ws.Range("R5").Value = 263.968836777512
$ws.Range("S5").Value = 0.182025935087863
$ws.Range("T5").Value = 0.1820259350878631
$ws.Range("I6").Value = 0.2200595722726403
$ws.Range("J6").Value = 0.2200595722726403
$ws.Range("O6").Value = 0.09421438109281059
$ws.Range("P6").Value = 0.09421438109281059
$ws.Range("S6").Value = 0.02073277640521543
$ws.Range("T6").Value = 0.02073277640521543
$ws.Range("I7").Value = 0.2200595722726403
$ws.Range("J7").Value = 0.2200595722726403
$ws.Range("O7").Value = 0.07861898758090437
$ws.Range("P7").Value = 0.07861898758090438
$ws.Range("S7").Value = 0.01730086077956184
$ws.Range("T7").Value = 0.01730086077956184
$ws.Range("H8").Value = 57.69206699999999
$ws.Range("I8").Value = 0.7155745972255653
$ws.Range("J8").Value = 0.7155745972255653
$ws.Range("M8").Value = 4.959409333333333
$ws.Range("N8").Value = 14.878228
$ws.Range("O8").Value = 0.8271666313262851
$ws.Range("P8").Value = 0.8271666313262852
$ws.Range("Q8").Value = 95.37285851303066
$ws.Range("R8").Value = 858.3557266172759
$ws.Range("S8").Value = 0.5918994290497341
$ws.Range("T8").Value = 0.5918994290497341
$ws.Range("H9").Value = 57.69206699999999
$ws.Range("I9").Value = 0.7155745972255653
$ws.Range("J9").Value = 0.7155745972255653
$ws.Range("O9").Value = 0.09421438109281059
$ws.Range("P9").Value = 0.09421438109281059
$ws.Range("R9").Value = 97.76682288434398
$ws.Range("S9").Value = 0.06741741780334386
$ws.Range("T9").Value = 0.06741741780334386
$ws.Range("H10").Value = 57.69206699999999
$ws.Range("I10").Value = 0.7155745972255653
$ws.Range("J10").Value = 0.7155745972255653
$ws.Range("O10").Value = 0.07861898758090437
$ws.Range("P10").Value = 0.07861898758090438
$ws.Range("Q10").Value = 9.064821155767332
$ws.Range("R10").Value = 81.58339040190599
$ws.Range("S10").Value = 0.05625775037248736
$ws.Range("T10").Value = 0.05625775037248737
